$d = $word.ActiveDocument

# Change 1: "So: the set is an array of Strings." -> append " No need for dictionary."
$p1 = $d.Paragraphs.Item(3)
$r1 = $p1.Range
$r1.Collapse(0)
$r1.MoveEnd(1, -1)
$r1.Collapse(0)
$r1.InsertAfter(" No need for dictionary.")

# Change 2a: append future-me note to the O(logn) paragraph (paragraph 8)
$p2 = $d.Paragraphs.Item(8)
$r2 = $p2.Range
$r2.Collapse(0)
$r2.MoveEnd(1, -1)
$r2.Collapse(0)
$r2.InsertAfter(" {Future me: It" + [char]0x2019 + "s actually O(nm), because there are 2 bases for O. It" + [char]0x2019 + "s like O(n^2), but not quite.}")

# Change 2b: add trailing space to the "Third:" paragraph (paragraph 10)
$p3 = $d.Paragraphs.Item(10)
$r3 = $p3.Range
$r3.Collapse(0)
$r3.MoveEnd(1, -1)
$r3.Collapse(0)
$r3.InsertAfter(" ")

# Change 3: add a new paragraph "{Future me: Nope.}" after "It's gotta be harder than it looks, then." (paragraph 16)
$p4 = $d.Paragraphs.Item(16)
$r4 = $p4.Range
$r4.Collapse(0)
$r4.InsertParagraphAfter()
$r4.Collapse(0)
$r4.MoveStart(1, 1)
$r4.Text = "{Future me: Nope.}"

# Change 4: add two new paragraphs at the end of the document
$last = $d.Paragraphs.Last
$rLast = $last.Range
$rLast.Collapse(0)
$rLast.InsertParagraphAfter()
$rLast.Collapse(0)
$rLast.MoveStart(1, 1)
$rLast.InsertParagraphAfter()
$rLast.Collapse(0)
$rLast.MoveStart(1, 1)
$rLast.Text = "{It is impossible for me, at the moment, to do this faster, simply because O(n) is the fastest way I can think of to manipulate an array. When doing it based on m, O(nm) is the fastest possible way. If there is one faster, though, I will gladly implement that instead.}"
